$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overview.Range("G2").Value = "2016-08-17 13:02:29"

$zhcn.Range("H2").Value = "2016-08-17 13:02:24"
$zhcn.Range("K2").Value = "2016-08-17 13:02:45"

$dede.Range("H2").Value = "2016-08-17 13:02:29"
$dede.Range("K2").Value = "2016-08-17 13:02:52"
